$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.874.72"
$ws.Range("E2").Value = "  -1.08%  "
$ws.Range("D3").Value = "3.536.51"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "611.47"
$ws.Range("E5").Value = "  +4.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "184.91"
$ws.Range("E6").Value = "  -0.67%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.620"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -0.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.645"
$ws.Range("E10").Value = "  -0.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.46"
$ws.Range("E11").Value = "  -1.28%  "
$ws.Range("E12").Value = "  -4.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.43"
$ws.Range("E13").Value = "  -1.21%  "
$ws.Range("D14").Value = "4.098.75"
$ws.Range("E14").Value = "  -1.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "621.44"
$ws.Range("E15").Value = "  +9.89%  "
$ws.Range("D16").Value = "69.905.38"
$ws.Range("E16").Value = "  -0.98%  "
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "12.62"
$ws.Range("E17").Value = "  +1.67%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.542.63"
$ws.Range("E18").Value = "  -1.76%  "
$ws.Range("E19").Value = "  -3.76%  "
$ws.Range("E20").Value = "  -0.38%  "
$ws.Range("E21").Value = "  -2.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.54"
$ws.Range("E22").Value = "  -0.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.72"
$ws.Range("E23").Value = "  +1.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "100.68"
$ws.Range("E24").Value = "  +5.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.91"
$ws.Range("E25").Value = "  +0.23%  "
$ws.Range("E26").Value = "  +1.49%  "
$ws.Range("E27").Value = "  -4.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.55"
$ws.Range("E28").Value = "  +4.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "32.77"
$ws.Range("E29").Value = "  +1.58%  "
$ws.Range("E30").Value = "  -4.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.21"
$ws.Range("E31").Value = "  -1.90%  "
$ws.Range("E32").Value = "  -0.67%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "63.44"
$ws.Range("E33").Value = "  -2.20%  "
$ws.Range("E34").Value = "  +14.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.25"
$ws.Range("E35").Value = "  -1.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "533.62"
$ws.Range("E36").Value = "  -5.22%  "
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.399"
$ws.Range("E38").Value = "  -4.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.07"
$ws.Range("E39").Value = "  -1.83%  "
$ws.Range("D40").Value = "0.0₃0779"
$ws.Range("E40").Value = "  -3.24%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "3.537.91"
$ws.Range("E41").Value = "  +4.68%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.56"
$ws.Range("E42").Value = "  +5.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.136"
$ws.Range("E43").Value = "  +1.47%  "
$ws.Range("E44").Value = "  +2.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.91"
$ws.Range("E45").Value = "  -1.77%  "
$ws.Range("E46").Value = "  +4.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.38"
$ws.Range("E47").Value = "  -4.96%  "
$ws.Range("E48").Value = "  -3.57%  "
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.42"
$ws.Range("E50").Value = "  -1.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "134.91"
$ws.Range("E51").Value = "  -2.10%  "
